$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.513.95"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").Value = "1.841.19"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.61%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9988"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07445"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2959"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").Value = "1.831.53"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.026"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6791"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009367"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.940"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("D18").Value = "29.458.65"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("D19").Value = "2.082.14"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.369"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9990"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1418"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.532"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06075"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.49%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("E31").Value = "  +2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.106"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.877"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7288"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.611"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").Value = "1.223.09"
$ws.Range("E39").Value = "  +1.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.993.40"
$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5077"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4066"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1141"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.58%  "
